# Swap bonferroni-correction significance markers with fdr-correction markers.
# For a selected set of cells that currently read "< .05", append an asterisk
# to read "< .05*" - mirroring the already-starred cells elsewhere in the table.
# Sheet "arousal" (sheet1) and sheet "valence" (sheet2) each get their own set
# of cell addresses updated.

$wb = $excel.ActiveWorkbook

$wsArousal = $wb.Worksheets.Item("arousal")
$wsValence = $wb.Worksheets.Item("valence")

$arousalCells = @(
  "C4", "C5", "C6", "C8", "C9", "C10", "C12", "H12", "H13", "C15", "H16",
  "C23", "H26", "H36", "H38", "H39", "H46", "C50", "H50", "C51", "C52",
  "C56", "C63", "C66"
)

$valenceCells = @(
  "H4", "H8", "H11", "H14", "H21", "H25", "H26", "H29", "H34", "H43",
  "H62", "H65", "H66"
)

foreach ($addr in $arousalCells) {
  $cell = $wsArousal.Range($addr)
  $cell.Value = $cell.Value() + "*"
}

foreach ($addr in $valenceCells) {
  $cell = $wsValence.Range($addr)
  $cell.Value = $cell.Value() + "*"
}
